$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.901.37'
$ws.Range("E2").Value = '  +0.47%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.694.37'
$ws.Range("E3").Value = '  +0.04%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '648.82'
$ws.Range("E5").Value = '  -4.14%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.80'
$ws.Range("E6").Value = '  +0.11%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.504'
$ws.Range("E8").Value = '  +1.70%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.147'
$ws.Range("E9").Value = '  -0.45%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.19'
$ws.Range("E10").Value = '  +0.92%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.446'
$ws.Range("E11").Value = '  +1.16%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000233'
$ws.Range("E12").Value = '  -0.29%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.318.26'
$ws.Range("E13").Value = '  +0.08%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.79'
$ws.Range("E14").Value = '  +0.87%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.686.68'
$ws.Range("E15").Value = '  +0.39%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.906.04'
$ws.Range("E16").Value = '  +0.59%  '

$ws.Range("E17").Value = '  +0.38%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '16.08'
$ws.Range("E18").Value = '  +0.34%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.52'
$ws.Range("E19").Value = '  +0.68%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.46'
$ws.Range("E20").Value = '  +6.70%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '473.02'
$ws.Range("E21").Value = '  +0.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.653'
$ws.Range("E22").Value = '  +0.40%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '80.17'
$ws.Range("E23").Value = '  -0.50%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.843.28'
$ws.Range("E24").Value = '  +0.13%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000128'
$ws.Range("E25").Value = '  +1.82%  '

$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.03'
$ws.Range("E27").Value = '  +1.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.18'
$ws.Range("E28").Value = '  +0.55%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.66'
$ws.Range("E29").Value = '  -1.57%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.72'
$ws.Range("E30").Value = '  -1.17%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.02'
$ws.Range("E31").Value = '  +0.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.57'
$ws.Range("E32").Value = '  -0.48%  '

$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  -0.24%  '

$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.93'
$ws.Range("E34").Value = '  -0.30%  '

$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.166'
$ws.Range("E35").Value = '  +2.50%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.692.25'
$ws.Range("E36").Value = '  +0.29%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.45'
$ws.Range("E37").Value = '  -0.22%  '

$ws.Range("E38").Value = '  -0.08%  '

$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '180.92'
$ws.Range("E39").Value = '  +7.62%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.26'
$ws.Range("E40").Value = '  +0.34%  '

$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.91'
$ws.Range("E41").Value = '  -5.18%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0909'
$ws.Range("E43").Value = '  +0.72%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.932'
$ws.Range("E44").Value = '  -1.15%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.84'
$ws.Range("E45").Value = '  +3.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '46.93'
$ws.Range("E46").Value = '  +0.84%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '29.18'
$ws.Range("E47").Value = '  +4.19%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000274'
$ws.Range("E48").Value = '  -1.69%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.27'
$ws.Range("E49").Value = '  -2.12%  '

$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.86'
$ws.Range("E50").Value = '  -0.40%  '

$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.06'
$ws.Range("E51").Value = '  -2.70%  '
